# Add a new column S (year 2022) to the "9.1.2" worksheet, mirroring the
# existing layout of column R (year 2021).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header - year value in row 3
$ws.Range("S3").Value = 2022

# Data rows 4-13 - numeric values copied from the source diff.
$ws.Range("S4").Value = 10444.200000000001
$ws.Range("S5").Value = 21.7
$ws.Range("S6").Value = 7361.6
$ws.Range("S7").Value = 143.1
$ws.Range("S8").Value = 844.2
# Row 9 intentionally left blank (S9 has no value in the source diff).
$ws.Range("S10").Value = "2 756,0"
$ws.Range("S11").Value = "1 013,8"
$ws.Range("S12").Value = "1 451,1"
$ws.Range("S13").Value = 273.39999999999998
$ws.Range("S14").Value = "-"
$ws.Range("S15").Value = 17.7

# Copy formatting from column R onto the new column S so the new cells pick
# up the same styles already used throughout the table.
$ws.Range("R1:R21").Copy()
$ws.Range("S1:S21").PasteSpecial(-4122)  # xlPasteFormats

# Set the active selection, matching the author's saved view.
$ws.Range("T3").Select()
